$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '87.972.30'
$ws.Range('E2').Value = '  +10.22%  '
$ws.Range('D3').Value = '3.329.97'
$ws.Range('E3').Value = '  +5.14%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.04'
$ws.Range('E5').Value = '  +5.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '637.40'
$ws.Range('E6').Value = '  +1.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.331'
$ws.Range('E7').Value = '  +22.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.613'
$ws.Range('E9').Value = '  +4.70%  '
$ws.Range('D10').Value = '3.328.17'
$ws.Range('E10').Value = '  +5.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.611'
$ws.Range('E11').Value = '  +3.79%  '
$ws.Range('E12').Value = '  +9.03%  '
$ws.Range('E13').Value = '  +2.06%  '
$ws.Range('D14').Value = '3.941.39'
$ws.Range('E14').Value = '  +5.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.35'
$ws.Range('E15').Value = '  +9.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.41'
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('D17').Value = '87.532.46'
$ws.Range('E17').Value = '  +10.08%  '
$ws.Range('D18').Value = '3.329.40'
$ws.Range('E18').Value = '  +5.96%  '
$ws.Range('B19').Value = 'SuiNetwork'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.21'
$ws.Range('E19').Value = '  +6.68%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.54'
$ws.Range('E20').Value = '  +2.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '448.04'
$ws.Range('E21').Value = '  +3.24%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.35'
$ws.Range('E23').Value = '  +3.61%  '
$ws.Range('E24').Value = '  +7.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.35'
$ws.Range('E25').Value = '  +14.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.29'
$ws.Range('E26').Value = '  +13.75%  '
$ws.Range('D27').Value = '3.511.70'
$ws.Range('E27').Value = '  +6.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '78.58'
$ws.Range('E28').Value = '  +3.78%  '
$ws.Range('E29').Value = '  +6.71%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  +53.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '599.85'
$ws.Range('E32').Value = '  +8.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.33'
$ws.Range('E33').Value = '  +4.77%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.54'
$ws.Range('E35').Value = '  +5.22%  '
$ws.Range('E36').Value = '  +3.51%  '
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.43'
$ws.Range('E38').Value = '  +1.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.59'
$ws.Range('E39').Value = '  +18.27%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.419'
$ws.Range('E40').Value = '  +3.73%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  +3.05%  '
$ws.Range('E43').Value = '  +14.58%  '
$ws.Range('E44').Value = '  +13.63%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '157.03'
$ws.Range('E46').Value = '  -3.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '188.90'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.21'
$ws.Range('E48').Value = '  +9.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.37'
$ws.Range('E49').Value = '  +6.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.787'
$ws.Range('E50').Value = '  +1.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.54'
$ws.Range('E51').Value = '  +8.71%  '
